# Applies crypto price/volume updates per commit "Updated cryptos list on Mon May 22 22:40:23 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.007.23"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.828.70"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'311.47"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "'0.4627"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'0.3699"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "'0.8765"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "'0.07879"
$ws.Range("E11").Value = "  +3.48%  "
$ws.Range("D12").Value = "'19.74"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "1.823.72"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "'5.336"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "'6.558"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'91.51"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "'1.007"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "'0.000008837"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "'14.78"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "27.030.89"
$ws.Range("D22").Value = "'5.102"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").Value = "'10.51"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").Value = "2.065.44"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").Value = "'152.26"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "'1.852"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'18.45"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("D28").Value = "'2.038"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "'5.108"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'115.70"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").Value = "'0.08886"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'2.964"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'0.7292"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "'4.437"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "'1.131"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.079"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'2.457"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("D38").Value = "'0.01943"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").Value = "'0.05218"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").Value = "'2.957"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Value = "'7.106"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").Value = "'0.5156"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").Value = "'0.1626"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "'8.155"
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").Value = "'0.4823"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").Value = "'10.13"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").Value = "'101.97"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "'1.625"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "'0.06207"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  +0.25%  "
